$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything shifts right by one.
$ws.Columns.Item(1).Insert()

# New column A: "Match ID" header in row 3, value 12 for data rows 4-18,
# and value 12 (default style) for the hidden summary row 19.
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

$ws.Range("A4:A18").Value = 12
$ws.Range("A4:A18").Font.Bold = $true

$ws.Range("A19").Value = 12
$ws.Rows.Item(19).AutoFit()

# Update the active selection to match the new target state.
$ws.Range("A3:A18").Select() | Out-Null
